$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-09-22 Monday"; new="2025-09-23 Tuesday"},
    @{old="516×7="; new="716×3="},
    @{old="425×2="; new="467×5="},
    @{old="838×3="; new="135×3="},
    @{old="824×4="; new="715×3="},
    @{old="531×7="; new="245×9="},
    @{old="154×7="; new="968×4="},
    @{old="730×8="; new="648×3="},
    @{old="857×8="; new="619×2="},
    @{old="320×7="; new="898×3="},
    @{old="807×9="; new="664×3="},
    @{old="850×5="; new="577×9="},
    @{old="427×8="; new="956×8="},
    @{old="584×4="; new="393×3="},
    @{old="529×6="; new="429×8="},
    @{old="760×8="; new="639×2="},
    @{old="935×9="; new="421×2="},
    @{old="416×6="; new="307×7="},
    @{old="568×6="; new="467×9="},
    @{old="731×3="; new="972×5="},
    @{old="259×7="; new="540×3="},
    @{old="599×8="; new="760×3="},
    @{old="241×2="; new="578×7="},
    @{old="559×3="; new="891×8="},
    @{old="175×8="; new="822×5="},
    @{old="498×7="; new="786×5="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
